$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.133.88'
$ws.Range('E2').Value = '  -3.51%  '
$ws.Range('D3').Value = '3.136.26'
$ws.Range('E3').Value = '  -4.87%  '
$ws.Range('E4').Value = '  +0.05%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.11'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  -5.84%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.78'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -4.89%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.135.79'
$ws.Range('E8').Value = '  -4.84%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.444'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('E10').Value = '  -7.33%  '
$ws.Range('E11').Value = '  -8.48%  '
$ws.Range('E12').Value = '  -6.44%  '
$ws.Range('D13').Value = '3.674.14'
$ws.Range('E13').Value = '  -4.89%  '
$ws.Range('E14').Value = '  -1.29%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.57'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  -4.53%  '
$ws.Range('D16').Value = '3.139.67'
$ws.Range('E16').Value = '  -4.85%  '
$ws.Range('D17').Value = '58.130.37'
$ws.Range('E17').Value = '  -3.56%  '
$ws.Range('E18').Value = '  -7.63%  '
$ws.Range('E19').Value = '  -4.84%  '
$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.08'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  -6.52%  '
$ws.Range('E21').Value = '  -7.81%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '344.11'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  -8.10%  '
$ws.Range('E23').Value = '  +0.13%  '
$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.47'
$ws.Range('D24').Style = $__style
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.508'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  -4.97%  '
$ws.Range('D26').Value = '3.266.86'
$ws.Range('E26').Value = '  -5.10%  '
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = '0.0₃0956'
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('E29').Value = '  +0.22%  '
$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.79'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  -4.92%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -8.65%  '
$ws.Range('E33').Value = '  -9.13%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.50'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -4.69%  '
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('E36').Value = '  -6.02%  '
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.00'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  -5.76%  '
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.22'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  -6.64%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.38'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  -9.36%  '
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0690'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  -4.90%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '24.62'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  -8.28%  '
$ws.Range('B42').Value = 'RenzoRestakedETH'
$ws.Range('C42').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D42').Value = '3.167.58'
$ws.Range('E42').Value = '  -4.84%  '
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.34'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('E44').Value = '  -7.85%  '
$ws.Range('E45').Value = '  -2.06%  '
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.90'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  -5.28%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  -8.36%  '
$ws.Range('D49').Value = '2.259.98'
$ws.Range('E49').Value = '  -3.78%  '
$ws.Range('E50').Value = '  -3.23%  '
$__style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.68'
$ws.Range('D51').Style = $__style
$ws.Range('E51').Value = '  -2.76%  '
